# Update "想去人数" (F column) values across sheets to reflect newly generated data.
$wb = $excel.ActiveWorkbook

function Set-FValue($SheetName, $Row, $Value) {
    $ws = $wb.Worksheets.Item($SheetName)
    $ws.Cells.Item($Row, 6).Value = $Value
}

# 展览 (sheet1)
Set-FValue "展览" 2  3418
Set-FValue "展览" 5  837
Set-FValue "展览" 7  274
Set-FValue "展览" 9  162
Set-FValue "展览" 10 636
Set-FValue "展览" 12 431
Set-FValue "展览" 13 68
Set-FValue "展览" 14 494
Set-FValue "展览" 15 335
Set-FValue "展览" 18 97
Set-FValue "展览" 19 185

# 演出 (sheet2)
Set-FValue "演出" 9 176

# 本地生活 (sheet3)
Set-FValue "本地生活" 2 6207
Set-FValue "本地生活" 3 740
Set-FValue "本地生活" 4 743
Set-FValue "本地生活" 5 1786
Set-FValue "本地生活" 6 117

# 全部类型 (sheet4)
Set-FValue "全部类型" 2  6207
Set-FValue "全部类型" 3  740
Set-FValue "全部类型" 4  743
Set-FValue "全部类型" 5  1786
Set-FValue "全部类型" 6  3418
Set-FValue "全部类型" 7  117
Set-FValue "全部类型" 10 837
Set-FValue "全部类型" 12 274
Set-FValue "全部类型" 17 162
Set-FValue "全部类型" 20 636
Set-FValue "全部类型" 24 431
Set-FValue "全部类型" 25 176
Set-FValue "全部类型" 26 68
Set-FValue "全部类型" 27 494
Set-FValue "全部类型" 29 335
Set-FValue "全部类型" 34 97
Set-FValue "全部类型" 40 185
